$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C13").Value = "Recherche bzgl. Frontend JavaScript SinglePageArchitecture"
$ws.Range("C17").Value = "Tests mit SinglePageArchitecture und Implementierung einer groben Struktur"
$ws.Range("C21").Value = "Ausbau Frontend für CRUD Operations"
$ws.Range("C25").Value = "Ausbau Frontend für CRUD Operations"
$ws.Range("C29").Value = "Fehlerbehebung bei Ungleichheiten in der Datenbankstruktur und im Frontend"
$ws.Range("C33").Value = "Fehlerbehebung bei Ungleichheiten in der Datenbankstruktur und im Frontend"

$ws.Range("C33").Select()
